# Append newly scraped Lancers listings as of 2025-10-20 01:23:53.
#
# The scraper always writes its freshest rows starting at row 11 (just
# below the previously-newest items) and pushes everything that was
# there before further down the sheet. Two new postings came in this
# run, inserted at two different points, and every row's "fetched at"
# timestamp is refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = '2025-10-20 01:23:53'

# ---------------------------------------------------------------
# 1) Refresh the "fetched at" timestamp for every already-known row
#    (rows 2-10 keep their content, only column A changes).
# ---------------------------------------------------------------
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# ---------------------------------------------------------------
# 2) Insert the first new posting at row 11. This pushes the old
#    row 11 ("高額成功報酬...") down to row 12, and the old row 12
#    ("人気調査...") down to row 13 (for now).
# ---------------------------------------------------------------
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = $newTimestamp
$ws.Cells.Item(11, 2).Value = 'サイトスピードが遅く サイトスピードを速くしたい ワードプレス'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5416402'
$ws.Cells.Item(11, 7).Value = 30
$ws.Cells.Item(11, 8).Value = '◇サイト'

# The row that used to be at 11 is now at 12; stamp its refreshed time.
$ws.Cells.Item(12, 1).Value = $newTimestamp

# ---------------------------------------------------------------
# 3) Insert the second new posting at row 13 (below the row that
#    used to be row 11). This pushes the old row 12 ("人気調査...",
#    currently sitting at row 13) down to row 14.
# ---------------------------------------------------------------
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = $newTimestamp
$ws.Cells.Item(13, 2).Value = '【急募】エクセルマクロの組み方を教えてください!'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5416433'
$ws.Cells.Item(13, 7).Value = 10

# The row that used to be row 12 is now at 14; stamp its refreshed time.
$ws.Cells.Item(14, 1).Value = $newTimestamp

# ---------------------------------------------------------------
# 4) Rebuild the hyperlinks on column F for every data row so each
#    URL cell links to its own address (row inserts do not shift the
#    hyperlink collection automatically).
# ---------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$lastRow = 14
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}

$wb.Save()
